$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 126.252295589447
$ws.Range("C2").Value = 2.532700610986595
$ws.Range("D2").Value = 1.115542268753052
$ws.Range("E2").Value = 0.1634631134554174
$ws.Range("J2").Value = 0.7787234042553192
$ws.Range("K2").Value = 0.723404255319149
$ws.Range("L2").Value = 0.7574468085106383
$ws.Range("M2").Value = 0.825531914893617
$ws.Range("N2").Value = 0.7905982905982906
$ws.Range("O2").Value = 0.7751409347154028
$ws.Range("P2").Value = 0.03400229674394581

$ws.Range("B3").Value = 246.872664642334
$ws.Range("C3").Value = 1.353709186814066
$ws.Range("D3").Value = 1.053472280502319
$ws.Range("E3").Value = 0.0424919047404017
$ws.Range("J3").Value = 0.774468085106383
$ws.Range("K3").Value = 0.7319148936170212
$ws.Range("L3").Value = 0.7617021276595745
$ws.Range("M3").Value = 0.8042553191489362
$ws.Range("N3").Value = 0.8076923076923077
$ws.Range("O3").Value = 0.7760065466448445
$ws.Range("P3").Value = 0.02811774434154455

$ws.Range("B4").Value = 489.2687978744507
$ws.Range("C4").Value = 2.088638246224441
$ws.Range("D4").Value = 1.214611005783081
$ws.Range("E4").Value = 0.1512280779726354
$ws.Range("J4").Value = 0.7914893617021277
$ws.Range("K4").Value = 0.7617021276595745
$ws.Range("L4").Value = 0.774468085106383
$ws.Range("M4").Value = 0.8127659574468085
$ws.Range("N4").Value = 0.7692307692307693
$ws.Range("O4").Value = 0.7819312602291326
$ws.Range("P4").Value = 0.0182666642565243
$ws.Range("Q4").Value = 7

$ws.Range("B5").Value = 126.0701298236847
$ws.Range("C5").Value = 0.7379490144250626
$ws.Range("D5").Value = 1.044546985626221
$ws.Range("E5").Value = 0.03617091231430056
$ws.Range("J5").Value = 0.7787234042553192
$ws.Range("K5").Value = 0.723404255319149
$ws.Range("L5").Value = 0.7574468085106383
$ws.Range("M5").Value = 0.825531914893617
$ws.Range("N5").Value = 0.7905982905982906
$ws.Range("O5").Value = 0.7751409347154028
$ws.Range("P5").Value = 0.03400229674394581

$ws.Range("B6").Value = 248.2222033977509
$ws.Range("C6").Value = 2.115320521666607
$ws.Range("D6").Value = 1.09271354675293
$ws.Range("E6").Value = 0.0547458299221872
$ws.Range("J6").Value = 0.774468085106383
$ws.Range("K6").Value = 0.7319148936170212
$ws.Range("L6").Value = 0.7617021276595745
$ws.Range("M6").Value = 0.8042553191489362
$ws.Range("N6").Value = 0.8076923076923077
$ws.Range("O6").Value = 0.7760065466448445
$ws.Range("P6").Value = 0.02811774434154455

$ws.Range("B7").Value = 491.3675094127655
$ws.Range("C7").Value = 1.481394093642365
$ws.Range("D7").Value = 1.147619724273682
$ws.Range("E7").Value = 0.07331906913002409
$ws.Range("J7").Value = 0.7914893617021277
$ws.Range("K7").Value = 0.7617021276595745
$ws.Range("L7").Value = 0.774468085106383
$ws.Range("M7").Value = 0.8127659574468085
$ws.Range("N7").Value = 0.7692307692307693
$ws.Range("O7").Value = 0.7819312602291326
$ws.Range("P7").Value = 0.0182666642565243
$ws.Range("Q7").Value = 7

$ws.Range("B8").Value = 126.6165143489838
$ws.Range("C8").Value = 1.070895491870624
$ws.Range("D8").Value = 1.061511564254761
$ws.Range("E8").Value = 0.06258126757202499
$ws.Range("J8").Value = 0.7787234042553192
$ws.Range("K8").Value = 0.723404255319149
$ws.Range("L8").Value = 0.7574468085106383
$ws.Range("M8").Value = 0.825531914893617
$ws.Range("N8").Value = 0.7905982905982906
$ws.Range("O8").Value = 0.7751409347154028
$ws.Range("P8").Value = 0.03400229674394581

$ws.Range("B9").Value = 246.0936964511871
$ws.Range("C9").Value = 1.194584778358449
$ws.Range("D9").Value = 1.04787392616272
$ws.Range("E9").Value = 0.02689566432830435
$ws.Range("J9").Value = 0.774468085106383
$ws.Range("K9").Value = 0.7319148936170212
$ws.Range("L9").Value = 0.7617021276595745
$ws.Range("M9").Value = 0.8042553191489362
$ws.Range("N9").Value = 0.8076923076923077
$ws.Range("O9").Value = 0.7760065466448445
$ws.Range("P9").Value = 0.02811774434154455

$ws.Range("B10").Value = 487.1323921203613
$ws.Range("C10").Value = 2.620914021985151
$ws.Range("D10").Value = 1.054905319213867
$ws.Range("E10").Value = 0.06393930322073368
$ws.Range("J10").Value = 0.7914893617021277
$ws.Range("K10").Value = 0.7617021276595745
$ws.Range("L10").Value = 0.774468085106383
$ws.Range("M10").Value = 0.8127659574468085
$ws.Range("N10").Value = 0.7692307692307693
$ws.Range("O10").Value = 0.7819312602291326
$ws.Range("P10").Value = 0.0182666642565243
$ws.Range("Q10").Value = 7

$ws.Range("B11").Value = 200.7952862262726
$ws.Range("C11").Value = 0.9305540677825398
$ws.Range("D11").Value = 1.079675674438477
$ws.Range("E11").Value = 0.07946847080706726
$ws.Range("J11").Value = 0.7957446808510639
$ws.Range("K11").Value = 0.7404255319148936
$ws.Range("L11").Value = 0.7531914893617021
$ws.Range("M11").Value = 0.8170212765957446
$ws.Range("N11").Value = 0.811965811965812
$ws.Range("O11").Value = 0.7836697581378432
$ws.Range("P11").Value = 0.03116980004243094
$ws.Range("Q11").Value = 1

$ws.Range("B12").Value = 393.3219874382019
$ws.Range("C12").Value = 2.638630151651159
$ws.Range("D12").Value = 1.080135345458984
$ws.Range("E12").Value = 0.08319124918889872
$ws.Range("J12").Value = 0.774468085106383
$ws.Range("K12").Value = 0.7659574468085106
$ws.Range("L12").Value = 0.7574468085106383
$ws.Range("M12").Value = 0.8042553191489362
$ws.Range("N12").Value = 0.811965811965812
$ws.Range("O12").Value = 0.782818694308056
$ws.Range("P12").Value = 0.02147952644971909

$ws.Range("B13").Value = 755.2623684406281
$ws.Range("C13").Value = 1.866882201860358
$ws.Range("D13").Value = 1.099753379821777
$ws.Range("E13").Value = 0.1151702403595287
$ws.Range("J13").Value = 0.7702127659574468
$ws.Range("K13").Value = 0.7574468085106383
$ws.Range("L13").Value = 0.7404255319148936
$ws.Range("M13").Value = 0.8
$ws.Range("N13").Value = 0.8076923076923077
$ws.Range("O13").Value = 0.7751554828150573
$ws.Range("P13").Value = 0.02537751797291475
$ws.Range("Q13").Value = 13

$ws.Range("B14").Value = 200.424781370163
$ws.Range("C14").Value = 0.9011002734419271
$ws.Range("D14").Value = 1.057651853561401
$ws.Range("E14").Value = 0.02752248455894615
$ws.Range("J14").Value = 0.7957446808510639
$ws.Range("K14").Value = 0.7404255319148936
$ws.Range("L14").Value = 0.7531914893617021
$ws.Range("M14").Value = 0.8170212765957446
$ws.Range("N14").Value = 0.811965811965812
$ws.Range("O14").Value = 0.7836697581378432
$ws.Range("P14").Value = 0.03116980004243094
$ws.Range("Q14").Value = 1

$ws.Range("B15").Value = 393.9644478797912
$ws.Range("C15").Value = 1.714323568784391
$ws.Range("D15").Value = 1.035970020294189
$ws.Range("E15").Value = 0.04371067858185509
$ws.Range("J15").Value = 0.774468085106383
$ws.Range("K15").Value = 0.7659574468085106
$ws.Range("L15").Value = 0.7574468085106383
$ws.Range("M15").Value = 0.8042553191489362
$ws.Range("N15").Value = 0.811965811965812
$ws.Range("O15").Value = 0.782818694308056
$ws.Range("P15").Value = 0.02147952644971909

$ws.Range("B16").Value = 753.3830567359925
$ws.Range("C16").Value = 4.839301410888211
$ws.Range("D16").Value = 1.095749568939209
$ws.Range("E16").Value = 0.1151560397422201
$ws.Range("J16").Value = 0.7702127659574468
$ws.Range("K16").Value = 0.7574468085106383
$ws.Range("L16").Value = 0.7404255319148936
$ws.Range("M16").Value = 0.8
$ws.Range("N16").Value = 0.8076923076923077
$ws.Range("O16").Value = 0.7751554828150573
$ws.Range("P16").Value = 0.02537751797291475
$ws.Range("Q16").Value = 13

$ws.Range("B17").Value = 200.1276122570038
$ws.Range("C17").Value = 0.8741474506619059
$ws.Range("D17").Value = 1.155200719833374
$ws.Range("E17").Value = 0.1730897323332244
$ws.Range("J17").Value = 0.7957446808510639
$ws.Range("K17").Value = 0.7404255319148936
$ws.Range("L17").Value = 0.7531914893617021
$ws.Range("M17").Value = 0.8170212765957446
$ws.Range("N17").Value = 0.811965811965812
$ws.Range("O17").Value = 0.7836697581378432
$ws.Range("P17").Value = 0.03116980004243094
$ws.Range("Q17").Value = 1

$ws.Range("B18").Value = 386.0691440582276
$ws.Range("C18").Value = 6.54838030377497
$ws.Range("D18").Value = 0.9614781856536865
$ws.Range("E18").Value = 0.07330760726739061
$ws.Range("J18").Value = 0.774468085106383
$ws.Range("K18").Value = 0.7659574468085106
$ws.Range("L18").Value = 0.7574468085106383
$ws.Range("M18").Value = 0.8042553191489362
$ws.Range("N18").Value = 0.811965811965812
$ws.Range("O18").Value = 0.782818694308056
$ws.Range("P18").Value = 0.02147952644971909

$ws.Range("B19").Value = 617.5517903327942
$ws.Range("C19").Value = 29.4161269021902
$ws.Range("D19").Value = 0.762101697921753
$ws.Range("E19").Value = 0.1764677928616137
$ws.Range("J19").Value = 0.7702127659574468
$ws.Range("K19").Value = 0.7574468085106383
$ws.Range("L19").Value = 0.7404255319148936
$ws.Range("M19").Value = 0.8
$ws.Range("N19").Value = 0.8076923076923077
$ws.Range("O19").Value = 0.7751554828150573
$ws.Range("P19").Value = 0.02537751797291475
$ws.Range("Q19").Value = 7
